# Update "想去人数" (attendance/interest count) figures in column F
# for both the "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 338
$ws1.Range("F4").Value = 2887
$ws1.Range("F5").Value = 70
$ws1.Range("F6").Value = 608

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 338
$ws4.Range("F6").Value = 2887
$ws4.Range("F7").Value = 70
$ws4.Range("F8").Value = 608
